$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the tracking number in I1
$ws.Range("I1").Value = 33085

# Mark "x" and fill in the bill details for Nicor Gas (row 8, no structural change)
$ws.Range("B8").Value = "x"
$ws.Range("D8").Value = 148.17

# Insert the 4 new vendor rows. They are written in the order the strings
# were originally entered (HUB, Salvi, Hernandez, ComCast) -- each one is
# inserted directly after its alphabetical predecessor in the list as it
# currently stands -- so that both the final alphabetical ordering of rows
# and the shared-string table creation order come out correct.

# After "Hinckley Springs" (currently row 6) -> HUB International Midwest Limited
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "HUB International Midwest Limited"
$ws.Range("G7").Value = "x"

# After "Pitney Bowes Purchase Power" (currently row 10) -> Salvi Salvi & Wifler
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Salvi Salvi & Wifler"

# After "Garlock Chicago Inc." (currently row 5) -> Hernandez Lawn Service
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Hernandez Lawn Service"
$ws.Range("G6").Value = "x"

# After "Blue Cross Blue Shield of Illinois" (currently row 2) -> ComCast Business
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "ComCast Business"
$ws.Range("B3").Value = "x"
$ws.Range("D3").Value = 471.04

# Update the selection to reflect the final state
$ws.Range("B13").Select()
